# The "DateJoined" row (row 4) currently stores real Excel date serials
# (B4:D4) formatted with a date number format. The fix re-enters those
# values as plain text ("2022-03-12", "2022-03-13", "2022-03-15") instead
# of dates, which is what the regression test now expects.
#
# Formatting the range as Text *before* typing the values is the standard
# way (in the Excel UI and via COM) to stop Excel from auto-parsing a
# date-shaped string into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4:D4").NumberFormat = "@"

$ws.Range("B4").Value = "2022-03-12"
$ws.Range("C4").Value = "2022-03-13"
$ws.Range("D4").Value = "2022-03-15"
